$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.911.74'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '1.874.26'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7407'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.37'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3146'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07151'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.69'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08412'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7506'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.402'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').Value = '1.862.15'
$ws.Range('E14').Value = '  -12.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.46'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.01%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '29.911.39'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.104'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.58'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.08'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007806'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.87%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').Value = '2.124.35'
$ws.Range('E22').Value = '  -6.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.983'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1555'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.294'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.45'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.90%  '
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.037'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.486'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.617'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.530'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.258'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05326'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.237'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7535'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9962'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.696'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01948'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.751'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4489'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('D42').Value = '1.112.59'
$ws.Range('E42').Value = '  +0.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.056'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.17'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8562'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.658'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.074'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.837'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.51%  '
$ws.Range('D51').Value = '2.022.70'
$ws.Range('E51').Value = '  -7.41%  '
